$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E hold numeric-looking text (prices / percentages) that must
# stay stored as text, matching the source data. Pre-format as Text ("@")
# before assigning so Excel does not auto-convert them to numbers.
$textCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D27", "E27", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "E46", "D47", "E47", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Plain text cells (coin names / links) -- assign directly.
$ws.Range("B4").Value = "LEO"
$ws.Range("C4").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("B5").Value = "HuobiToken"
$ws.Range("C5").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("B6").Value = "Cronos"
$ws.Range("C6").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"

# Numeric-looking text cells -- assign now that the cells are Text-formatted.
$ws.Range("D2").Value = "309.27"
$ws.Range("E2").Value = "1.32%"
$ws.Range("D3").Value = "39.22"
$ws.Range("E3").Value = "9.65%"
$ws.Range("D4").Value = "3.595"
$ws.Range("E4").Value = "-3.60%"
$ws.Range("D5").Value = "5.113"
$ws.Range("E5").Value = "1.49%"
$ws.Range("D6").Value = "0.08142"
$ws.Range("E6").Value = "2.20%"
$ws.Range("D7").Value = "1.981"
$ws.Range("E7").Value = "3.72%"
$ws.Range("D8").Value = "7.917"
$ws.Range("E8").Value = "1.89%"
$ws.Range("D9").Value = "0.9332"
$ws.Range("E9").Value = "1.43%"
$ws.Range("D10").Value = "0.1407"
$ws.Range("E10").Value = "7.98%"
$ws.Range("D11").Value = "0.1941"
$ws.Range("E11").Value = "1.56%"
$ws.Range("D12").Value = "0.09148"
$ws.Range("E12").Value = "0.57%"
$ws.Range("D13").Value = "0.03485"
$ws.Range("E13").Value = "1.44%"
$ws.Range("D14").Value = "0.09823"
$ws.Range("E14").Value = "-0.07%"
$ws.Range("D15").Value = "0.001422"
$ws.Range("E15").Value = "1.45%"
$ws.Range("D16").Value = "0.005839"
$ws.Range("E16").Value = "-5.15%"
$ws.Range("D17").Value = "4.195"
$ws.Range("E17").Value = "1.48%"
$ws.Range("E18").Value = "2.77%"
$ws.Range("D19").Value = "0.3447"
$ws.Range("E19").Value = "0.06%"
$ws.Range("D20").Value = "0.1323"
$ws.Range("E20").Value = "0.97%"
$ws.Range("D21").Value = "4.822"
$ws.Range("E21").Value = "-6.57%"
$ws.Range("D22").Value = "0.2469"
$ws.Range("E22").Value = "5.09%"
$ws.Range("D23").Value = "0.04459"
$ws.Range("E23").Value = "1.14%"
$ws.Range("D24").Value = "0.001240"
$ws.Range("E24").Value = "0.65%"
$ws.Range("D25").Value = "0.004861"
$ws.Range("E25").Value = "4.98%"
$ws.Range("D27").Value = "0.0001304"
$ws.Range("E27").Value = "4.33%"
$ws.Range("D39").Value = "0.02121"
$ws.Range("E39").Value = "9.09%"
$ws.Range("D40").Value = "0.05124"
$ws.Range("E40").Value = "-2.87%"
$ws.Range("D41").Value = "0.007481"
$ws.Range("E41").Value = "-1.41%"
$ws.Range("D42").Value = "0.009946"
$ws.Range("E42").Value = "-1.70%"
$ws.Range("D43").Value = "0.1363"
$ws.Range("E43").Value = "0.82%"
$ws.Range("D44").Value = "0.002136"
$ws.Range("E44").Value = "-1.08%"
$ws.Range("D45").Value = "0.01013"
$ws.Range("E45").Value = "1.86%"
$ws.Range("E46").Value = "2.00%"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").Value = "0.30%"
$ws.Range("E48").Value = "-0.63%"
$ws.Range("D49").Value = "0.001603"
$ws.Range("E49").Value = "-3.25%"
$ws.Range("D50").Value = "0.00002106"
$ws.Range("E50").Value = "0.30%"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").Value = "0.30%"

# Restore the default "Normal" style on the text-forced cells so only the
# number format needed to keep them textual is applied (no stray styling).
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
